# Testdata_Signup.xlsx edit: clear stale sample values from the
# Valid_Testcases sheet (B11:B19), drop the now-orphaned hyperlinks that
# used to sit on B16/B18/B19, and update the saved sheet
# selections/active tab to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$wsInvalid = $wb.Worksheets.Item("Invalid_Testcases")
$wsValid   = $wb.Worksheets.Item("Valid_Testcases")

# --- Clear out the leftover lastname/email-ish values in column B -------
$wsValid.Range("B11:B19").ClearContents()

# --- Remove the hyperlinks that used to be attached to B16/B18/B19 ------
$targets = @('$B$16', '$B$18', '$B$19')
$toRemove = @()
foreach ($hl in $wsValid.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($targets -contains $addr) {
        $toRemove += $hl
    }
}
# Delete back-to-front so earlier references in $toRemove stay valid as the
# live Hyperlinks collection re-indexes after each removal.
for ($i = $toRemove.Count - 1; $i -ge 0; $i--) {
    $toRemove[$i].Delete()
}

# --- Update sheet selections / active sheet -----------------------------
# Valid_Testcases is no longer the active tab; its selection moved to I47.
$wsValid.Activate() | Out-Null
$wsValid.Range("I47").Select() | Out-Null

# Invalid_Testcases becomes the active/selected tab, selection at C46.
$wsInvalid.Activate() | Out-Null
$wsInvalid.Range("C46").Select() | Out-Null
